$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so that numeric-looking
# strings (e.g. "0.9989") are not auto-converted to numbers, matching the
# original inline-string cell content. Restore the default "Normal" style
# afterwards so no stray cell-style attribute is introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$priceUpdates = @{
    'D2' = '28.242.32'
    'D3' = '1.803.39'
    'D5' = '339.04'
    'D6' = '0.9989'
    'D7' = '0.4887'
    'D8' = '0.3715'
    'D9' = '45.59'
    'D10' = '0.07697'
    'D11' = '1.150'
    'D12' = '22.68'
    'D13' = '0.9999'
    'D14' = '6.328'
    'D15' = '7.333'
    'D16' = '1.798.08'
    'D17' = '0.00001099'
    'D18' = '0.06719'
    'D19' = '82.44'
    'D20' = '0.9998'
    'D21' = '17.45'
    'D22' = '6.432'
    'D23' = '28.236.78'
    'D24' = '12.03'
    'D25' = '2.410'
    'D26' = '20.90'
    'D27' = '2.418'
    'D28' = '151.19'
    'D29' = '2.004.09'
    'D30' = '134.59'
    'D31' = '1.273'
    'D32' = '4.036'
    'D33' = '0.09733'
    'D34' = '5.972'
    'D35' = '0.02383'
    'D36' = '12.23'
    'D37' = '0.2208'
    'D38' = '0.06353'
    'D39' = '0.6717'
    'D40' = '5.254'
    'D41' = '1.487'
    'D43' = '8.149'
    'D44' = '14.22'
    'D45' = '0.9989'
    'D46' = '0.6185'
    'D47' = '3.872'
    'D48' = '129.45'
    'D49' = '2.057'
    'D50' = '1.174'
    'D51' = '0.07111'
}

foreach ($cell in $priceUpdates.Keys) {
    $ws.Range($cell).Value = $priceUpdates[$cell]
}

$priceRange.Style = "Normal"

# Remaining updates: Volume(1h) percentages (column E, always textual because
# of the "%" sign) and the Coin name / Link swap for rows 37-38 (column B/C).
$otherUpdates = @{
    'E2' = '  +0.91%  '
    'E3' = '  +2.40%  '
    'E4' = '  -0.09%  '
    'E5' = '  +0.55%  '
    'E6' = '  +0.01%  '
    'E7' = '  +29.37%  '
    'E8' = '  +10.74%  '
    'E9' = '  -0.31%  '
    'E10' = '  +7.05%  '
    'E11' = '  +2.51%  '
    'E12' = '  +1.57%  '
    'E13' = '  -0.12%  '
    'E14' = '  +2.17%  '
    'E15' = '  +1.97%  '
    'E16' = '  +2.22%  '
    'E17' = '  +4.32%  '
    'E18' = '  +2.14%  '
    'E19' = '  +2.59%  '
    'E20' = '  +0.04%  '
    'E21' = '  +2.72%  '
    'E22' = '  +2.38%  '
    'E23' = '  +0.86%  '
    'E24' = '  +2.57%  '
    'E25' = '  +1.63%  '
    'E26' = '  +5.22%  '
    'E27' = '  +3.38%  '
    'E28' = '  -0.95%  '
    'E29' = '  +2.23%  '
    'E30' = '  +2.06%  '
    'E31' = '  +1.56%  '
    'E32' = '  +0.52%  '
    'E33' = '  +10.66%  '
    'E34' = '  +2.88%  '
    'E35' = '  +1.73%  '
    'E36' = '  -0.38%  '
    'B37' = 'Algorand'
    'C37' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'E37' = '  +4.42%  '
    'B38' = 'Hedera'
    'C38' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'E38' = '  +2.90%  '
    'E39' = '  +1.67%  '
    'E40' = '  +1.97%  '
    'E41' = '  +2.69%  '
    'E42' = '  +0.77%  '
    'E43' = '  +1.63%  '
    'E44' = '  +3.08%  '
    'E45' = '  +0.01%  '
    'E46' = '  +2.07%  '
    'E47' = '  +1.27%  '
    'E48' = '  -0.29%  '
    'E49' = '  +2.23%  '
    'E50' = '  -0.70%  '
    'E51' = '  -0.72%  '
}

foreach ($cell in $otherUpdates.Keys) {
    $ws.Range($cell).Value = $otherUpdates[$cell]
}
